$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

Set-TextValue 'D2' '43.945.25'
Set-TextValue 'E2' '  -0.68%  '
Set-TextValue 'D3' '2.234.27'
Set-TextValue 'E4' '  +0.20%  '
Set-TextValue 'D5' '312.30'
Set-TextValue 'E5' '  -2.48%  '
Set-TextValue 'D6' '98.61'
Set-TextValue 'E6' '  -5.34%  '
Set-TextValue 'E7' '  -3.90%  '
Set-TextValue 'E8' '  +0.23%  '
Set-TextValue 'D9' '0.534'
Set-TextValue 'E9' '  -7.45%  '
Set-TextValue 'D10' '36.13'
Set-TextValue 'E10' '  -5.75%  '
Set-TextValue 'D11' '0.0819'
Set-TextValue 'E11' '  -3.28%  '
Set-TextValue 'D12' '7.35'
Set-TextValue 'E12' '  -7.26%  '
Set-TextValue 'E13' '  -3.24%  '
Set-TextValue 'D14' '2.571.54'
Set-TextValue 'E14' '  -2.00%  '
Set-TextValue 'D15' '2.274.32'
Set-TextValue 'E15' '  -0.66%  '
Set-TextValue 'D16' '0.836'
Set-TextValue 'E16' '  -5.07%  '
Set-TextValue 'D17' '14.09'
Set-TextValue 'E17' '  -3.68%  '
Set-TextValue 'D18' '43.796.38'
Set-TextValue 'E18' '  -0.77%  '
Set-TextValue 'D19' '12.88'
Set-TextValue 'E19' '  -12.24%  '
Set-TextValue 'D20' '0.0₃0960'
Set-TextValue 'E20' '  -4.12%  '
Set-TextValue 'E21' '  -5.37%  '
Set-TextValue 'D22' '64.90'
Set-TextValue 'E22' '  -2.35%  '
Set-TextValue 'E23' '  -7.16%  '
Set-TextValue 'D24' '232.75'
Set-TextValue 'E24' '  -3.07%  '
Set-TextValue 'D25' '2.03'
Set-TextValue 'E25' '  -9.47%  '
Set-TextValue 'E26' '  +0.23%  '
Set-TextValue 'D27' '10.11'
Set-TextValue 'E27' '  -1.47%  '
Set-TextValue 'D28' '2.16'
Set-TextValue 'E28' '  -1.69%  '
Set-TextValue 'D29' '36.64'
Set-TextValue 'E29' '  -7.67%  '
Set-TextValue 'D30' '5.94'
Set-TextValue 'E30' '  -9.25%  '
Set-TextValue 'D31' '157.84'
Set-TextValue 'E31' '  -2.32%  '
Set-TextValue 'D32' '19.89'
Set-TextValue 'E32' '  -3.74%  '
Set-TextValue 'D33' '0.0828'
Set-TextValue 'E33' '  -7.02%  '
Set-TextValue 'E34' '  -1.50%  '
Set-TextValue 'E35' '  -6.61%  '
Set-TextValue 'E36' '  +1.37%  '
Set-TextValue 'E37' '  -6.96%  '
Set-TextValue 'E38' '  -4.15%  '
Set-TextValue 'D39' '15.56'
Set-TextValue 'E39' '  -3.02%  '
Set-TextValue 'D40' '3.60'
Set-TextValue 'E40' '  -9.89%  '
Set-TextValue 'D41' '4.03'
Set-TextValue 'E41' '  -10.98%  '
Set-TextValue 'E42' '  -7.21%  '
Set-TextValue 'E43' '  +0.08%  '
Set-TextValue 'D44' '1.711.68'
Set-TextValue 'E44' '  -5.82%  '
Set-TextValue 'E45' '  -8.12%  '
Set-TextValue 'D46' '80.09'
Set-TextValue 'E46' '  -9.00%  '
Set-TextValue 'E47' '  -0.91%  '
Set-TextValue 'B48' 'ordi'
Set-TextValue 'C48' 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
Set-TextValue 'D48' '72.97'
Set-TextValue 'E48' '  -5.13%  '
Set-TextValue 'B49' 'THORChain'
Set-TextValue 'C49' 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
Set-TextValue 'D49' '5.08'
Set-TextValue 'E49' '  -6.71%  '
Set-TextValue 'D50' '101.36'
Set-TextValue 'E50' '  -3.48%  '
Set-TextValue 'D51' '56.17'
Set-TextValue 'E51' '  -6.72%  '
